$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.755.95'
$ws.Range('E2').Value = '  +10.20%  '

$ws.Range('D3').Value = '3.338.46'
$ws.Range('E3').Value = '  +4.56%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('D5').Value = '220.40'
$ws.Range('E5').Value = '  +4.77%  '

$ws.Range('D6').Value = '653.55'
$ws.Range('E6').Value = '  +2.38%  '

$ws.Range('D7').Value = '0.360'
$ws.Range('E7').Value = '  +24.59%  '

$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').Value = '0.612'
$ws.Range('E9').Value = '  +3.35%  '

$ws.Range('D10').Value = '3.336.86'
$ws.Range('E10').Value = '  +4.56%  '

$ws.Range('D11').Value = '0.597'
$ws.Range('E11').Value = '  +0.86%  '

$ws.Range('D12').Value = '0.0000273'
$ws.Range('E12').Value = '  +2.81%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '35.79'
$ws.Range('E13').Value = '  +11.43%  '

$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.168'
$ws.Range('E14').Value = '  +1.76%  '

$ws.Range('D15').Value = '3.932.44'
$ws.Range('E15').Value = '  +4.02%  '

$ws.Range('D16').Value = '5.45'
$ws.Range('E16').Value = '  +2.09%  '

$ws.Range('D17').Value = '88.339.91'
$ws.Range('E17').Value = '  +9.75%  '

$ws.Range('D18').Value = '3.306.42'
$ws.Range('E18').Value = '  +3.64%  '

$ws.Range('D19').Value = '14.74'
$ws.Range('E19').Value = '  +2.41%  '

$ws.Range('D20').Value = '3.15'
$ws.Range('E20').Value = '  +1.40%  '

$ws.Range('D21').Value = '465.93'
$ws.Range('E21').Value = '  +4.16%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '9.23'
$ws.Range('E22').Value = '  -0.41%  '

$ws.Range('B23').Value = 'Polkadot'
$ws.Range('C23').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D23').Value = '5.57'
$ws.Range('E23').Value = '  +6.10%  '

$ws.Range('B24').Value = 'NEARProtocol'
$ws.Range('C24').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D24').Value = '5.56'
$ws.Range('E24').Value = '  +15.60%  '

$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').Value = '12.91'
$ws.Range('E25').Value = '  +17.31%  '

$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.487.77'
$ws.Range('E26').Value = '  +3.79%  '

$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value = '79.37'
$ws.Range('E27').Value = '  +2.43%  '

$ws.Range('B28').Value = 'Cronos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D28').Value = '0.219'
$ws.Range('E28').Value = '  +79.33%  '

$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0000128'
$ws.Range('E29').Value = '  +3.67%  '

$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '613.72'
$ws.Range('E31').Value = '  +7.55%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '9.49'
$ws.Range('E32').Value = '  +4.26%  '

$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.61'
$ws.Range('E33').Value = '  +10.12%  '

$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.01%  '

$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  +2.38%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.150'
$ws.Range('E36').Value = '  -1.53%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '6.91'
$ws.Range('E37').Value = '  +20.49%  '

$ws.Range('D38').Value = '23.82'
$ws.Range('E38').Value = '  +3.92%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '2.19'
$ws.Range('E39').Value = '  +15.91%  '

$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '0.423'
$ws.Range('E40').Value = '  +2.89%  '

$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D41').Value = '21.79'
$ws.Range('E41').Value = '  +4.72%  '

$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.19%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '3.08'
$ws.Range('E43').Value = '  +11.42%  '

$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('D45').Value = '190.25'
$ws.Range('E45').Value = '  +0.77%  '

$ws.Range('D46').Value = '155.77'
$ws.Range('E46').Value = '  -2.12%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '47.20'
$ws.Range('E47').Value = '  +10.11%  '

$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').Value = '1.41'
$ws.Range('E48').Value = '  +7.57%  '

$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.797'
$ws.Range('E49').Value = '  +2.75%  '

$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '4.47'
$ws.Range('E50').Value = '  +4.31%  '

$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '0.665'
$ws.Range('E51').Value = '  +4.90%  '
